# Updated cryptos list on Wed Mar 13 03:40:07 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be stored as text so that numeric-looking strings
    # (e.g. "537.35") are not auto-converted into floating point numbers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "72.005.03"
$ws.Range("E2").Value = "  +0.25%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "4.037.64"
$ws.Range("E3").Value = "  -0.30%  "

# Row 5 - BNB
Set-TextValue "D5" "537.35"
$ws.Range("E5").Value = "  +0.77%  "

# Row 6 - Solana
Set-TextValue "D6" "149.09"
$ws.Range("E6").Value = "  -3.22%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "4.032.63"
$ws.Range("E7").Value = "  -0.30%  "

# Row 8 - XRP
Set-TextValue "D8" "0.696"
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - USDC (D unchanged)
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.751"
$ws.Range("E10").Value = "  -1.54%  "

# Row 11 - Dogecoin (D unchanged)
$ws.Range("E11").Value = "  -2.71%  "

# Row 12 - Avalanche
Set-TextValue "D12" "53.51"
$ws.Range("E12").Value = "  +8.20%  "

# Row 13 - ShibaInu (D unchanged)
$ws.Range("E13").Value = "  -1.66%  "

# Row 14 - Polkadot
Set-TextValue "D14" "10.89"
$ws.Range("E14").Value = "  -1.18%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.685.79"
$ws.Range("E15").Value = "  -0.11%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "4.050.68"
$ws.Range("E16").Value = "  +0.18%  "

# Row 17 - Uniswap
Set-TextValue "D17" "14.30"
$ws.Range("E17").Value = "  -1.55%  "

# Row 18 - Chainlink
Set-TextValue "D18" "20.71"
$ws.Range("E18").Value = "  -2.00%  "

# Row 19 - Polygon (D unchanged)
$ws.Range("E19").Value = "  -2.06%  "

# Row 20 - TRON (D unchanged)
$ws.Range("E20").Value = "  -1.26%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "72.019.79"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "436.90"
$ws.Range("E22").Value = "  -0.28%  "

# Row 23 - Litecoin
Set-TextValue "D23" "98.10"
$ws.Range("E23").Value = "  -1.76%  "

# Row 24 - ImmutableX
Set-TextValue "D24" "3.51"
$ws.Range("E24").Value = "  -6.29%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "4.26"
$ws.Range("E25").Value = "  +0.22%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "14.62"
$ws.Range("E26").Value = "  -1.62%  "

# Row 27 - Toncoin
Set-TextValue "D27" "4.39"
$ws.Range("E27").Value = "  +23.24%  "

# Row 28 - RenderToken
Set-TextValue "D28" "11.27"
$ws.Range("E28").Value = "  -1.70%  "

# Row 29 - Filecoin
Set-TextValue "D29" "10.72"
$ws.Range("E29").Value = "  -2.05%  "

# Row 30 - LEO
Set-TextValue "D30" "5.95"
$ws.Range("E30").Value = "  +2.18%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "37.07"
$ws.Range("E31").Value = "  -0.90%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "8.33"
$ws.Range("E32").Value = "  +22.83%  "

# Row 33 - Hedera (D unchanged)
$ws.Range("E33").Value = "  +1.87%  "

# Row 34 - InjectiveProtocol
Set-TextValue "D34" "50.17"
$ws.Range("E34").Value = "  +15.34%  "

# Row 35 - Cosmos
Set-TextValue "D35" "13.56"
$ws.Range("E35").Value = "  -1.17%  "

# Row 36 - Bittensor
Set-TextValue "D36" "683.49"
$ws.Range("E36").Value = "  +1.03%  "

# Row 37 - OKB
Set-TextValue "D37" "66.96"
$ws.Range("E37").Value = "  +0.38%  "

# Row 38 - TheGraph
Set-TextValue "D38" "0.460"
$ws.Range("E38").Value = "  +5.52%  "

# Row 39 - PEPE
$ws.Range("D39").Value = "0.0₃0861"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40 - WEMIXToken
Set-TextValue "D40" "3.43"
$ws.Range("E40").Value = "  +8.14%  "

# Row 41 - Kaspa (D unchanged)
$ws.Range("E41").Value = "  -6.29%  "

# Row 42 - ThetaToken
Set-TextValue "D42" "3.39"
$ws.Range("E42").Value = "  -1.36%  "

# Row 43 - THORChain
Set-TextValue "D43" "11.19"
$ws.Range("E43").Value = "  +16.62%  "

# Row 44 - Dai
Set-TextValue "D44" "1.00"
$ws.Range("E44").Value = "  +0.18%  "

# Row 45 - was VeChain, now FirstDigitalUSD
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D45" "1.00"
$ws.Range("E45").Value = "  +0.23%  "

# Row 46 - was FirstDigitalUSD, now VeChain
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0493"
$ws.Range("E46").Value = "  -1.43%  "

# Row 47 - Stellar (D unchanged)
$ws.Range("E47").Value = "  -1.17%  "

# Row 48 - Fetch.AI (D unchanged)
$ws.Range("E48").Value = "  -2.97%  "

# Row 49 - Stacks (D unchanged)
$ws.Range("E49").Value = "  +1.24%  "

# Row 50 - ApeXProtocol
Set-TextValue "D50" "3.30"
$ws.Range("E50").Value = "  -3.15%  "

# Row 51 - was FLOKI, now Maker
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.844.55"
$ws.Range("E51").Value = "  +9.15%  "
